$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 60.961853333015583
$ws.Range("C2").Value = 45.976771617899388
$ws.Range("D2").Value = 59.700097501423443
$ws.Range("E2").Value = 48.79964776451186

$ws.Range("B3").Value = 61.48547715904084
$ws.Range("C3").Value = 42.602687163599157
$ws.Range("D3").Value = 64.43982137233499
$ws.Range("E3").Value = 43.767603607131896

$ws.Range("B1:E3").Select()
